# Chandra rows cleaned up (count 13)
# Remove the row whose Question/Answer pair is the short
# "Okay, great. Thanks for taking my question." /
# "Go back and follow up. We can follow up." exchange - this row is
# row 42 in the original sheet (Asker "Jay McCanless"). Deleting it
# shifts all subsequent rows up by one and the now-unused shared
# strings are dropped automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$ws.Rows.Item(42).Delete()

# Leave the selection where the author ended up after the cleanup.
$ws.Range("C44").Select()
